# Weekly update: two new Pomelo price records (2022-08-08, serial 44811) are
# inserted at the top of the data block (rows 15-16), pushing all of the
# existing history down by two rows (old row 15 -> new row 17, ... old row 40
# -> new row 42). Row 14 (and the header in row 1) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 15 - this shifts rows
# 15:40 down to 17:42, carrying their values/formatting along, and extends
# the sheet's used range (dimension) accordingly.
$ws.Rows("15:16").Insert()

# Fill in the two new records.
$newRows = @(
    @{ Row = 15; A = 6; B = "Mercado Mayorista Lo Valledor de Santiago"; C = "Metropolitana";
       D = 44811; E = 13; F = "Fruta"; G = 100102; H = "Cítricos"; I = 100102006; J = "Pomelo";
       K = "Start Ruby"; L = "Especial"; M = 12; N = 170000; O = 170000; P = 170000;
       Q = "`$/bins (350 kilos)"; R = "Provincia de Limarí"; S = 486; T = 350 },
    @{ Row = 16; A = 6; B = "Mercado Mayorista Lo Valledor de Santiago"; C = "Metropolitana";
       D = 44811; E = 13; F = "Fruta"; G = 100102; H = "Cítricos"; I = 100102006; J = "Pomelo";
       K = "Start Ruby"; L = "Primera"; M = 25; N = 150000; O = 150000; P = 150000;
       Q = "`$/bins (350 kilos)"; R = "Provincia de Limarí"; S = 429; T = 350 }
)

foreach ($rec in $newRows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value  = $rec.A
    $ws.Cells.Item($r, 2).Value  = $rec.B
    $ws.Cells.Item($r, 3).Value  = $rec.C
    $ws.Cells.Item($r, 4).Value  = $rec.D
    $ws.Cells.Item($r, 5).Value  = $rec.E
    $ws.Cells.Item($r, 6).Value  = $rec.F
    $ws.Cells.Item($r, 7).Value  = $rec.G
    $ws.Cells.Item($r, 8).Value  = $rec.H
    $ws.Cells.Item($r, 9).Value  = $rec.I
    $ws.Cells.Item($r, 10).Value = $rec.J
    $ws.Cells.Item($r, 11).Value = $rec.K
    $ws.Cells.Item($r, 12).Value = $rec.L
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 15).Value = $rec.O
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
    $ws.Cells.Item($r, 18).Value = $rec.R
    $ws.Cells.Item($r, 19).Value = $rec.S
    $ws.Cells.Item($r, 20).Value = $rec.T
}
